# Event Service Release.xlsx -- add two new Use Case / Expected Behavior
# blocks (for the LessonModal write-operations -> LessonCalendar update
# use case) above the existing "user clicks on a single lesson" block,
# and renumber/shift everything below it down accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-LabelCell($row, $text) {
    # "Use Case:" / "Expected Behavior:" style labels -> bold Calibri 11
    $c = $ws.Cells.Item($row, 1)
    $c.Value = $text
    $c.Font.Name = "Calibri"
    $c.Font.Size = 11
    $c.Font.Bold = $true
}

function Set-BodyCell($row, $text) {
    # plain body text under a label -> regular Calibri 11
    $c = $ws.Cells.Item($row, 1)
    $c.Value = $text
    $c.Font.Name = "Calibri"
    $c.Font.Size = 11
    $c.Font.Bold = $false
}

# --- Remove the two stray formatting-only rows that previously lived
# between row 72 and the first "Use Case:" block (old A76 wrap-text
# placeholder, old A78 placeholder). Deleting (rather than merely
# clearing) removes the row entirely and shifts everything below up,
# which keeps later row-math simple.
$ws.Rows(78).Delete()
$ws.Rows(76).Delete()

# At this point the old "Use Case:" block that used to start at row 82
# now starts at row 80 (82 - 2). We want it to end up starting at row
# 87, so make room by inserting 7 blank rows right before it.
for ($i = 0; $i -lt 7; $i++) {
    $ws.Rows(80).Insert()
}

# Rows 74-86 are now free. Populate the two new use-case blocks.
Set-LabelCell 74 "Use Case: "
Set-BodyCell  75 "user performs write operations on a lesson (Edit/Delete) in LessonModal component"
Set-LabelCell 76 "Expected Behavior:"
Set-BodyCell  77 "LessonCalendar Element updates when LessonModal Closes"

# Row 78 stays blank but with a slightly taller custom height (spacer row).
$ws.Rows(78).RowHeight = 16

Set-LabelCell 79 "Use Case: "
Set-LabelCell 81 "Expected Behavior:"
Set-LabelCell 83 "Use Case: "
Set-LabelCell 85 "Expected Behavior:"

# Rows 87 onward already hold the original "user clicks on a single
# lesson..." use case block through the end of the sheet (shifted down
# by +5 automatically by the row delete/insert operations above), so
# no further edits are required there.

# Update the active selection to match the new layout.
$ws.Range("A77").Select() | Out-Null

Write-Host "Edit complete"
